# Update NATMI L1cam-Egfr TPM-based edge weight metrics (recomputed with new TPM values)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.713252999999999
$ws.Range("H2").Value = 20.139759
$ws.Range("I2").Value = 0.3101840064655811
$ws.Range("J2").Value = 0.3231642354899327
$ws.Range("M2").Value = 0.6598136666666666
$ws.Range("N2").Value = 1.979441
$ws.Range("O2").Value = 0.007704735356083927
$ws.Range("P2").Value = 0.008484678519943686
$ws.Range("Q2").Value = 4.429496077190999
$ws.Range("R2").Value = 39.86546469471899
$ws.Range("S2").Value = 0.002389885681507129
$ws.Range("T2").Value = 0.002741944647275455

# Row 3
$ws.Range("G3").Value = 6.713252999999999
$ws.Range("H3").Value = 20.139759
$ws.Range("I3").Value = 0.3101840064655811
$ws.Range("J3").Value = 0.3231642354899327
$ws.Range("O3").Value = 0.7130079175842846
$ws.Range("P3").Value = 0.7851850431306702
$ws.Range("Q3").Value = 409.9122978249779
$ws.Range("R3").Value = 3689.210680424801
$ws.Range("S3").Value = 0.2211636525179743
$ws.Range("T3").Value = 0.2537437241814529

# Row 4
$ws.Range("G4").Value = 6.713252999999999
$ws.Range("H4").Value = 20.139759
$ws.Range("I4").Value = 0.3101840064655811
$ws.Range("J4").Value = 0.3231642354899327
$ws.Range("K4").Value = 2.0
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1563486666666667
$ws.Range("N4").Value = 0.469046
$ws.Range("O4").Value = 0.001825704984300993
$ws.Range("P4").Value = 0.002010519394650058
$ws.Range("Q4").Value = 1.049608155546
$ws.Range("R4").Value = 9.446473399913998
$ws.Range("S4").Value = 0.000566304486654663
$ws.Range("T4").Value = 0.0006497279631097684

# Row 5
$ws.Range("G5").Value = 6.713252999999999
$ws.Range("H5").Value = 20.139759
$ws.Range("I5").Value = 0.3101840064655811
$ws.Range("J5").Value = 0.3231642354899327
$ws.Range("M5").Value = 23.6163295
$ws.Range("N5").Value = 47.232659
$ws.Range("O5").Value = 0.2757711427815902
$ws.Range("P5").Value = 0.2024581319964196
$ws.Range("Q5").Value = 158.5423948648635
$ws.Range("R5").Value = 951.2543691891809
$ws.Range("S5").Value = 0.08553979793558546
$ws.Range("T5").Value = 0.06542722744534282

# Row 6
$ws.Range("G6").Value = 6.713252999999999
$ws.Range("H6").Value = 20.139759
$ws.Range("I6").Value = 0.3101840064655811
$ws.Range("J6").Value = 0.3231642354899327
$ws.Range("K6").Value = 2.0
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.14477
$ws.Range("N6").Value = 0.43431
$ws.Range("O6").Value = 0.00169049929374041
$ws.Range("P6").Value = 0.001861626958316384
$ws.Range("Q6").Value = 0.97187763681
$ws.Range("R6").Value = 8.746898731289999
$ws.Range("S6").Value = 0.0005243658438596356
$ws.Range("T6").Value = 0.0006016112527517631

# Row 7
$ws.Range("I7").Value = 0.0154484264788496
$ws.Range("J7").Value = 0.01609489473505086
$ws.Range("M7").Value = 0.6598136666666666
$ws.Range("N7").Value = 1.979441
$ws.Range("O7").Value = 0.007704735356083927
$ws.Range("P7").Value = 0.008484678519943686
$ws.Range("Q7").Value = 0.2206069399468889
$ws.Range("R7").Value = 1.985462459522
$ws.Range("S7").Value = 0.0001190260376874556
$ws.Range("T7").Value = 0.0001365600076392407

# Row 8
$ws.Range("I8").Value = 0.0154484264788496
$ws.Range("J8").Value = 0.01609489473505086
$ws.Range("O8").Value = 0.7130079175842846
$ws.Range("P8").Value = 0.7851850431306702
$ws.Range("S8").Value = 0.01101485039363847
$ws.Range("T8").Value = 0.0126374706167245

# Row 9
$ws.Range("I9").Value = 0.0154484264788496
$ws.Range("J9").Value = 0.01609489473505086
$ws.Range("K9").Value = 2.0
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1563486666666667
$ws.Range("N9").Value = 0.469046
$ws.Range("O9").Value = 0.001825704984300993
$ws.Range("P9").Value = 0.002010519394650058
$ws.Range("Q9").Value = 0.05227475977022222
$ws.Range("R9").Value = 0.470472837932
$ws.Range("S9").Value = 0.00002820426922204315
$ws.Range("T9").Value = 0.00003235909801967086

# Row 10
$ws.Range("I10").Value = 0.0154484264788496
$ws.Range("J10").Value = 0.01609489473505086
$ws.Range("M10").Value = 23.6163295
$ws.Range("N10").Value = 47.232659
$ws.Range("O10").Value = 0.2757711427815902
$ws.Range("P10").Value = 0.2024581319964196
$ws.Range("Q10").Value = 7.896056791446333
$ws.Range("R10").Value = 47.376340748678
$ws.Range("S10").Value = 0.00426023022424973
$ws.Range("T10").Value = 0.003258542322737405

# Row 11
$ws.Range("I11").Value = 0.0154484264788496
$ws.Range("J11").Value = 0.01609489473505086
$ws.Range("K11").Value = 2.0
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.14477
$ws.Range("N11").Value = 0.43431
$ws.Range("O11").Value = 0.00169049929374041
$ws.Range("P11").Value = 0.001861626958316384
$ws.Range("Q11").Value = 0.04840346344666667
$ws.Range("R11").Value = 0.43563117102
$ws.Range("S11").Value = 0.00002611555405189589
$ws.Range("T11").Value = 0.00002996268993003512

# Row 12
$ws.Range("G12").Value = 6.661784666666667
$ws.Range("H12").Value = 19.985354
$ws.Range("I12").Value = 0.3078059262949933
$ws.Range("J12").Value = 0.3206866401135023
$ws.Range("M12").Value = 0.6598136666666666
$ws.Range("N12").Value = 1.979441
$ws.Range("O12").Value = 0.007704735356083927
$ws.Range("P12").Value = 0.008484678519943686
$ws.Range("Q12").Value = 4.395536567457111
$ws.Range("R12").Value = 39.55982910711401
$ws.Range("S12").Value = 0.002371563203137198
$ws.Range("T12").Value = 0.002720923047003944

# Row 13
$ws.Range("G13").Value = 6.661784666666667
$ws.Range("H13").Value = 19.985354
$ws.Range("I13").Value = 0.3078059262949933
$ws.Range("J13").Value = 0.3206866401135023
$ws.Range("O13").Value = 0.7130079175842846
$ws.Range("P13").Value = 0.7851850431306702
$ws.Range("Q13").Value = 406.7696331910236
$ws.Range("R13").Value = 3660.926698719212
$ws.Range("S13").Value = 0.219468062527695
$ws.Range("T13").Value = 0.25179835334895

# Row 14
$ws.Range("G14").Value = 6.661784666666667
$ws.Range("H14").Value = 19.985354
$ws.Range("I14").Value = 0.3078059262949933
$ws.Range("J14").Value = 0.3206866401135023
$ws.Range("K14").Value = 2.0
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1563486666666667
$ws.Range("N14").Value = 0.469046
$ws.Range("O14").Value = 0.001825704984300993
$ws.Range("P14").Value = 0.002010519394650058
$ws.Range("Q14").Value = 1.041561150253778
$ws.Range("R14").Value = 9.374050352284
$ws.Range("S14").Value = 0.0005619628138341533
$ws.Range("T14").Value = 0.0006447467095533597

# Row 15
$ws.Range("G15").Value = 6.661784666666667
$ws.Range("H15").Value = 19.985354
$ws.Range("I15").Value = 0.3078059262949933
$ws.Range("J15").Value = 0.3206866401135023
$ws.Range("M15").Value = 23.6163295
$ws.Range("N15").Value = 47.232659
$ws.Range("O15").Value = 0.2757711427815902
$ws.Range("P15").Value = 0.2024581319964196
$ws.Range("Q15").Value = 157.3269017460477
$ws.Range("R15").Value = 943.961410476286
$ws.Range("S15").Value = 0.08488399204931621
$ws.Range("T15").Value = 0.06492561811358775

# Row 16
$ws.Range("G16").Value = 6.661784666666667
$ws.Range("H16").Value = 19.985354
$ws.Range("I16").Value = 0.3078059262949933
$ws.Range("J16").Value = 0.3206866401135023
$ws.Range("K16").Value = 2.0
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.14477
$ws.Range("N16").Value = 0.43431
$ws.Range("O16").Value = 0.00169049929374041
$ws.Range("P16").Value = 0.001861626958316384
$ws.Range("Q16").Value = 0.9644265661933334
$ws.Range("R16").Value = 8.67983909574
$ws.Range("S16").Value = 0.0005203457010107989
$ws.Range("T16").Value = 0.0005969988944072003

# Row 17
$ws.Range("G17").Value = 2.607918
$ws.Range("H17").Value = 5.215835999999999
$ws.Range("I17").Value = 0.1204981331366039
$ws.Range("J17").Value = 0.08369373503331734
$ws.Range("M17").Value = 0.6598136666666666
$ws.Range("N17").Value = 1.979441
$ws.Range("O17").Value = 0.007704735356083927
$ws.Range("P17").Value = 0.008484678519943686
$ws.Range("Q17").Value = 1.720739937946
$ws.Range("R17").Value = 10.324439627676
$ws.Range("S17").Value = 0.0009284062267197002
$ws.Range("T17").Value = 0.000710114435891046

# Row 18
$ws.Range("G18").Value = 2.607918
$ws.Range("H18").Value = 5.215835999999999
$ws.Range("I18").Value = 0.1204981331366039
$ws.Range("J18").Value = 0.08369373503331734
$ws.Range("O18").Value = 0.7130079175842846
$ws.Range("P18").Value = 0.7851850431306702
$ws.Range("Q18").Value = 159.239888608268
$ws.Range("R18").Value = 955.4393316496079
$ws.Range("S18").Value = 0.08591612298052383
$ws.Range("T18").Value = 0.06571506895190216

# Row 19
$ws.Range("G19").Value = 2.607918
$ws.Range("H19").Value = 5.215835999999999
$ws.Range("I19").Value = 0.1204981331366039
$ws.Range("J19").Value = 0.08369373503331734
$ws.Range("K19").Value = 2.0
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.1563486666666667
$ws.Range("N19").Value = 0.469046
$ws.Range("O19").Value = 0.001825704984300993
$ws.Range("P19").Value = 0.002010519394650058
$ws.Range("Q19").Value = 0.4077445020759999
$ws.Range("R19").Value = 2.446467012456
$ws.Range("S19").Value = 0.0002199940422664624
$ws.Range("T19").Value = 0.0001682678774951876

# Row 20
$ws.Range("G20").Value = 2.607918
$ws.Range("H20").Value = 5.215835999999999
$ws.Range("I20").Value = 0.1204981331366039
$ws.Range("J20").Value = 0.08369373503331734
$ws.Range("M20").Value = 23.6163295
$ws.Range("N20").Value = 47.232659
$ws.Range("O20").Value = 0.2757711427815902
$ws.Range("P20").Value = 0.2024581319964196
$ws.Range("Q20").Value = 61.58945079698099
$ws.Range("R20").Value = 246.357803187924
$ws.Range("S20").Value = 0.03322990787812945
$ws.Range("T20").Value = 0.01694447725464873

# Row 21
$ws.Range("G21").Value = 2.607918
$ws.Range("H21").Value = 5.215835999999999
$ws.Range("I21").Value = 0.1204981331366039
$ws.Range("J21").Value = 0.08369373503331734
$ws.Range("K21").Value = 2.0
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.14477
$ws.Range("N21").Value = 0.43431
$ws.Range("O21").Value = 0.00169049929374041
$ws.Range("P21").Value = 0.001861626958316384
$ws.Range("Q21").Value = 0.37754828886
$ws.Range("R21").Value = 2.26528973316
$ws.Range("S21").Value = 0.0002037020089644668
$ws.Range("T21").Value = 0.000155806513380212

# Row 22
$ws.Range("G22").Value = 5.325505333333333
$ws.Range("H22").Value = 15.976516
$ws.Range("I22").Value = 0.2460635076239721
$ws.Range("J22").Value = 0.2563604946281968
$ws.Range("M22").Value = 0.6598136666666666
$ws.Range("N22").Value = 1.979441
$ws.Range("O22").Value = 0.007704735356083927
$ws.Range("P22").Value = 0.008484678519943686
$ws.Range("Q22").Value = 3.513841200839555
$ws.Range("R22").Value = 31.624570807556
$ws.Range("S22").Value = 0.001895854207032445
$ws.Range("T22").Value = 0.002175136382134

# Row 23
$ws.Range("G23").Value = 5.325505333333333
$ws.Range("H23").Value = 15.976516
$ws.Range("I23").Value = 0.2460635076239721
$ws.Range("J23").Value = 0.2563604946281968
$ws.Range("O23").Value = 0.7130079175842846
$ws.Range("P23").Value = 0.7851850431306702
$ws.Range("Q23").Value = 325.1762041838497
$ws.Range("R23").Value = 2926.585837654648
$ws.Range("S23").Value = 0.1754452291644531
$ws.Range("T23").Value = 0.2012904260316406

# Row 24
$ws.Range("G24").Value = 5.325505333333333
$ws.Range("H24").Value = 15.976516
$ws.Range("I24").Value = 0.2460635076239721
$ws.Range("J24").Value = 0.2563604946281968
$ws.Range("K24").Value = 2.0
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.1563486666666667
$ws.Range("N24").Value = 0.469046
$ws.Range("O24").Value = 0.001825704984300993
$ws.Range("P24").Value = 0.002010519394650058
$ws.Range("Q24").Value = 0.8326356581928889
$ws.Range("R24").Value = 7.493720923735999
$ws.Range("S24").Value = 0.0004492393723236713
$ws.Range("T24").Value = 0.0005154177464720718

# Row 25
$ws.Range("G25").Value = 5.325505333333333
$ws.Range("H25").Value = 15.976516
$ws.Range("I25").Value = 0.2460635076239721
$ws.Range("J25").Value = 0.2563604946281968
$ws.Range("M25").Value = 23.6163295
$ws.Range("N25").Value = 47.232659
$ws.Range("O25").Value = 0.2757711427815902
$ws.Range("P25").Value = 0.2024581319964196
$ws.Range("Q25").Value = 125.7688887060073
$ws.Range("R25").Value = 754.613332236044
$ws.Range("S25").Value = 0.0678572146943093
$ws.Range("T25").Value = 0.05190226686010288

# Row 26
$ws.Range("G26").Value = 5.325505333333333
$ws.Range("H26").Value = 15.976516
$ws.Range("I26").Value = 0.2460635076239721
$ws.Range("J26").Value = 0.2563604946281968
$ws.Range("K26").Value = 2.0
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.14477
$ws.Range("N26").Value = 0.43431
$ws.Range("O26").Value = 0.00169049929374041
$ws.Range("P26").Value = 0.001861626958316384
$ws.Range("Q26").Value = 0.7709734071066667
$ws.Range("R26").Value = 6.93876066396
$ws.Range("S26").Value = 0.0004159701858536128
$ws.Range("T26").Value = 0.0004772476078471738
